$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above current row 5, shifting the existing row 5
# (Rukh Lviv vs Ch. Odesa) down to row 7.
$ws.Rows.Item(5).Resize(2).Insert()

# --- Updated odds for existing row 3 (Arda vs Septemvri Sofia) ---
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 3.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("AG3").Value = 251

# --- Updated odds for existing row 4 (Penang vs Sabah) ---
$ws.Range("G4").Value = 2.95
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 2.05
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 2.27
$ws.Range("L4").Value = 2.57
$ws.Range("W4").Value = 10.25
$ws.Range("X4").Value = 14.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 29
$ws.Range("AA4").Value = 18.5
$ws.Range("AB4").Value = 21
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 9.75
$ws.Range("AJ4").Value = 7.5
$ws.Range("AK4").Value = 16
$ws.Range("AL4").Value = 12.5
$ws.Range("AM4").Value = 17
$ws.Range("AN4").Value = 5.2
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 19.5
$ws.Range("AQ4").Value = 65
$ws.Range("AR4").Value = 80
$ws.Range("AS4").Value = 200
$ws.Range("AX4").Value = 4.2
$ws.Range("AY4").Value = 10.25
$ws.Range("AZ4").Value = 16
$ws.Range("BA4").Value = 35
$ws.Range("BB4").Value = 55

# --- New row 5: Kedah vs PDRM FC ---
$row5 = @('42feT1O9','25/10/2024','10:00','MALAYSIA - SUPER LEAGUE','Kedah','PDRM FC',1.7,3.65,4.05,2.22,2.25,4.25,1.03,10,1.2,4.05,1.62,2.02,1.31,3.23,1.64,2.2,7,7.6,7,11.5,10.5,17.5,12.5,6.4,11.25,40,250,11.75,20,11.25,50,27,27,3.7,8.25,15.5,26,50,175,3.05,6.9,55,51,6.1,22,25,110,120,250,51)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# --- New row 6: Kelantan DNFC vs Sri Pahang ---
$row6 = @('QHr4RugM','25/10/2024','10:00','MALAYSIA - SUPER LEAGUE','Kelantan DNFC','Sri Pahang',4.9,3.85,1.55,4.9,2.3,2.05,1.02,10.5,1.19,4.19,1.6,2.07,1.29,3.32,1.68,2.13,13.5,26,13,70,35,32,13,6.8,12,45,250,7,7,6.8,10,9.75,17,6.9,28,29,150,175,350,3.1,7.1,55,51,3.5,7.3,15,22,45,175,51)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}

# --- Row 7 (previously row 5, Rukh Lviv vs Ch. Odesa) updated odds ---
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3.35
$ws.Range("J7").Value = 2.35
$ws.Range("K7").Value = 1.98
$ws.Range("N7").Value = 7.05
$ws.Range("Q7").Value = 2.35
$ws.Range("R7").Value = 1.47
$ws.Range("S7").Value = 1.5
$ws.Range("T7").Value = 2.25
$ws.Range("U7").Value = 2.25
$ws.Range("V7").Value = 1.5
$ws.Range("X7").Value = 6.5
$ws.Range("Z7").Value = 12.5
$ws.Range("AC7").Value = 6.8
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 24
$ws.Range("AI7").Value = 25
$ws.Range("AN7").Value = 3.25
$ws.Range("AO7").Value = 8.5
$ws.Range("AU7").Value = 9
